$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the Price column (D) keeps text formatting so numeric-looking
# strings (e.g. "1.001") are not auto-converted to numbers.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "28.273.12"
$ws.Range("E2").Value = "  +2.60%  "

$ws.Range("D3").Value = "1.870.35"
$ws.Range("E3").Value = "  +1.39%  "

$ws.Range("E4").Value = "  -0.28%  "

$ws.Range("D5").Value = "336.63"
$ws.Range("E5").Value = "  +0.83%  "

$ws.Range("D6").Value = "1.001"
$ws.Range("E6").Value = "  -0.39%  "

$ws.Range("D7").Value = "0.4695"
$ws.Range("E7").Value = "  +1.12%  "

$ws.Range("D8").Value = "0.3919"
$ws.Range("E8").Value = "  +1.73%  "

$ws.Range("D9").Value = "47.36"
$ws.Range("E9").Value = "  +2.79%  "

$ws.Range("D10").Value = "0.07984"
$ws.Range("E10").Value = "  +0.76%  "

$ws.Range("D11").Value = "1.004"
$ws.Range("E11").Value = "  +0.90%  "

$ws.Range("D12").Value = "21.68"
$ws.Range("E12").Value = "  +0.74%  "

$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").Value = "1.884.52"
$ws.Range("E13").Value = "  +2.01%  "

$ws.Range("B14").Value = "Polkadot"
$ws.Range("C14").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D14").Value = "5.982"
$ws.Range("E14").Value = "  +0.80%  "

$ws.Range("D15").Value = "7.267"
$ws.Range("E15").Value = "  +2.12%  "

$ws.Range("D16").Value = "91.20"
$ws.Range("E16").Value = "  +2.45%  "

$ws.Range("D17").Value = "1.003"
$ws.Range("E17").Value = "  -0.40%  "

$ws.Range("D18").Value = "0.00001041"
$ws.Range("E18").Value = "  +0.49%  "

$ws.Range("D19").Value = "0.06584"
$ws.Range("E19").Value = "  -1.38%  "

$ws.Range("E20").Value = "  +3.40%  "

$ws.Range("D21").Value = "1.000"
$ws.Range("E21").Value = "  -0.44%  "

$ws.Range("D22").Value = "28.279.68"
$ws.Range("E22").Value = "  +2.64%  "

$ws.Range("D23").Value = "5.441"
$ws.Range("E23").Value = "  +0.96%  "

$ws.Range("D24").Value = "11.03"
$ws.Range("E24").Value = "  +1.13%  "

$ws.Range("D25").Value = "2.297"
$ws.Range("E25").Value = "  -0.69%  "

$ws.Range("D26").Value = "2.099.21"
$ws.Range("E26").Value = "  +1.57%  "

$ws.Range("D27").Value = "159.00"
$ws.Range("E27").Value = "  +0.06%  "

$ws.Range("D28").Value = "19.84"
$ws.Range("E28").Value = "  +1.66%  "

$ws.Range("D29").Value = "2.144"
$ws.Range("E29").Value = "  +1.89%  "

$ws.Range("D30").Value = "5.500"
$ws.Range("E30").Value = "  +1.61%  "

$ws.Range("D31").Value = "119.79"
$ws.Range("E31").Value = "  -0.15%  "

$ws.Range("D32").Value = "0.9764"
$ws.Range("E32").Value = "  -0.30%  "

$ws.Range("D33").Value = "0.09475"
$ws.Range("E33").Value = "  +0.73%  "

$ws.Range("D34").Value = "3.580"
$ws.Range("E34").Value = "  -0.40%  "

$ws.Range("D35").Value = "5.350"
$ws.Range("E35").Value = "  +1.21%  "

$ws.Range("D36").Value = "1.370"
$ws.Range("E36").Value = "  +1.95%  "

$ws.Range("D37").Value = "0.02262"
$ws.Range("E37").Value = "  +1.51%  "

$ws.Range("D38").Value = "0.06088"
$ws.Range("E38").Value = "  +0.88%  "

$ws.Range("D39").Value = "8.424"
$ws.Range("E39").Value = "  +1.37%  "

$ws.Range("D40").Value = "1.174"
$ws.Range("E40").Value = "  -0.90%  "

$ws.Range("D41").Value = "0.5950"
$ws.Range("E41").Value = "  +1.01%  "

$ws.Range("D42").Value = "0.9995"
$ws.Range("E42").Value = "  -0.41%  "

$ws.Range("D43").Value = "0.1879"
$ws.Range("E43").Value = "  +0.59%  "

$ws.Range("D44").Value = "10.37"
$ws.Range("E44").Value = "  +0.52%  "

$ws.Range("D45").Value = "1.304"
$ws.Range("E45").Value = "  +5.06%  "

$ws.Range("D46").Value = "0.5607"
$ws.Range("E46").Value = "  +0.33%  "

$ws.Range("D47").Value = "12.13"
$ws.Range("E47").Value = "  -0.16%  "

$ws.Range("D48").Value = "1.965"
$ws.Range("E48").Value = "  +3.05%  "

$ws.Range("D49").Value = "0.06890"
$ws.Range("E49").Value = "  +2.85%  "

$ws.Range("D50").Value = "110.69"
$ws.Range("E50").Value = "  -0.40%  "

$ws.Range("D51").Value = "1.978"
$ws.Range("E51").Value = "  +11.12%  "
